# Fix projection detail creation not working properly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update raw poll figures in row 2 (projection inputs)
$ws.Range("A2").Value = 30
$ws.Range("B2").Value = 32
$ws.Range("C2").Value = 14
$ws.Range("D2").Value = 7
$ws.Range("I2").Value = 8

# Update raw figures in row 10 (secondary projection inputs)
$ws.Range("B10").Value = 52

# Update the selected cell to match where the user left off
$ws.Range("N8").Select()

# Reflect the window's new screen position (as recorded in the saved view state)
$win = $excel.ActiveWindow
$win.Left = 12540
$win.Top = 4665
